$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A width: 50 -> 44 (ColumnWidth uses a slightly narrower scale
#     than the <col width> stored in the xml, so compensate by ~0.83) ---
$ws.Columns("A").ColumnWidth = 43.17

# --- Make room for a 5th "Bad Drivers" data row: push the old Totals row
#     (and everything below it) down by one, then fill the new row 7 ---
$ws.Rows(7).Insert()

# --- Bad Drivers data rows (A3:D7) ---
$ws.Range("A3").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.110.1.1"
$ws.Range("B3").Value = 2
$ws.Range("C3").Value = 412
$ws.Range("D3").Value = 98.2

$ws.Range("A4").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 23.40.0.4"
$ws.Range("B4").Value = 8
$ws.Range("C4").Value = 923
$ws.Range("D4").Value = 98.4

$ws.Range("A5").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.160.3.2"
$ws.Range("B5").Value = 3
$ws.Range("C5").Value = 614
$ws.Range("D5").Value = 98.6

$ws.Range("A6").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 23.60.1.2"
$ws.Range("B6").Value = 2
$ws.Range("C6").Value = 106
$ws.Range("D6").Value = 98.7

# New 5th data row - copy formatting from the row above, then set values/text
$ws.Range("A6:D6").Copy($ws.Range("A7"))
$ws.Range("A7").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.240.0.6"
$ws.Range("B7").Value = 11
$ws.Range("C7").Value = 1185
$ws.Range("D7").Value = 98.8

# --- Totals row, now row 8 ---
$ws.Range("B8").Value = 26
$ws.Range("C8").Value = 3240

# --- Remove the first "Good Drivers" data row (old row 15, "20.50.0.5"),
#     which after the insert above now sits at row 16 ---
$ws.Rows(16).Delete()

# --- Remove the trailing "Good Drivers" data rows (old rows 22-32) ---
$ws.Range("A22:A32").EntireRow.Delete()

# --- Update the remaining "Good Drivers" data rows (now rows 16-21) ---
$ws.Range("A16").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.60.2.1"
$ws.Range("B16").Value = 56018
$ws.Range("D16").Value = 100

$ws.Range("A17").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.50.1.1"
$ws.Range("B17").Value = 34244

$ws.Range("B18").Value = 442178

$ws.Range("A19").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.80.0.9"
$ws.Range("B19").Value = 77849
$ws.Range("D19").Value = 99.9

$ws.Range("A20").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.110.3.2"
$ws.Range("B20").Value = 59673

$ws.Range("A21").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.70.0.6"
$ws.Range("B21").Value = 113652

# The "Driver Vintage" column holds yyyy-mm-dd text, not real dates - format
# the cells as Text first so Excel doesn't reinterpret the literal strings.
$ws.Range("E19:E21").NumberFormat = "@"
$ws.Range("E19").Value = "2021-08-18"
$ws.Range("E20").Value = "2020-08-05"
$ws.Range("E21").Value = "2019-12-14"

# --- Keep the sheet's used range (and therefore <dimension>) extending out
#     to column J / row 26, matching the report's fixed layout area, without
#     leaving any visible formatting behind ---
$ws.Range("J26").Borders.Item(9).LineStyle = 1
$ws.Range("J26").Borders.Item(9).LineStyle = -4142
